# This script applies the "2025-09-23" attendance column to the grades
# export sheet, and fixes the highlight coloring so that:
#   - half-day values (0.5) are highlighted yellow
#   - full-day values (1)   are highlighted red
#   - the stray red highlight that had been left on the 2025-09-18
#     column (I5, I7, I22) is removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 10092543   # RGB(255,255,153) encoded as OLE BGR for Interior.Color
$red    = 10066431   # RGB(255,153,153) encoded as OLE BGR for Interior.Color

# ---------------------------------------------------------------------
# 1) Remove the (incorrect) red highlight that used to sit on the
#    2025-09-18 column for these three students, leaving the values
#    untouched.
# ---------------------------------------------------------------------
$ws.Range("I5").Style = "Normal"

$ws.Range("I7").Style = "Normal"

$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = "1"
$ws.Range("I22").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Add the new 2025-09-23 attendance column (J), copying the header
#    formatting from the existing header cells.
# ---------------------------------------------------------------------
$ws.Range("J1").NumberFormat = "@"
$ws.Range("J1").Value = "2025-09-23"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Default attendance value for the new column is "0" (present), for
# every student row.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.NumberFormat = "@"
    $cell.Value = "0"
    $cell.Style = "Normal"
}

# Half-day absences (0.5) -> yellow highlight
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "0.5"
$ws.Range("J2").Interior.Color = $yellow

$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "0.5"
$ws.Range("J2").Copy()
$ws.Range("J15").PasteSpecial(-4122)

# Full-day absences (1) -> red highlight
$ws.Range("J27").NumberFormat = "@"
$ws.Range("J27").Value = "1"
$ws.Range("J27").Interior.Color = $red

$ws.Range("J28").Value = 1
$ws.Range("J27").Copy()
$ws.Range("J28").PasteSpecial(-4122)

$excel.CutCopyMode = 0
